$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the sheet (data runs from row 2 to the last row).
$lastRow = $ws.UsedRange.Rows.Count - 1

# Column C holds the "Förändrad" (changed) date, stored as a serial date value
# of 45188 for every data row (2..lastRow). Bump it by one day to 45189.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45189
